$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.42
$ws.Range("C2").Value = 0.38

$ws.Range("B3").Value = 0.27
$ws.Range("C3").Value = 0.1

$ws.Range("B4").Value = 0.1
$ws.Range("C4").Value = 0.31

$ws.Range("B5").Value = 0.11
$ws.Range("C5").Value = 0.15

$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

$ws.Range("B7").Value = 0.09
$ws.Range("C7").Value = 0.05
